# "Generate Report for Archive"
# - Flip the localization status shown on the Overview sheet (per-locale
#   columns) and on each locale's own status column from
#   "Ready for handoff" to "In Translation".
# - Narrow the now-shorter status columns to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns for each locale (zh-cn = E, de-de = F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Columns shrink now that "In Translation" is shorter than "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766

# --- Per-locale sheets: their own "Status" column (column C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.576851254417766

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.576851254417766
